# Auto-generated edit script: apply scheduled-runner market-data refresh
# to the Anima Profits workbook (columns H-N per leve row, per sheet).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1365
$ws.Range("I98").Value = 1365
$ws.Range("K98").Value = 1365
$ws.Range("M98").Value = 133
$ws.Range("H103").Value = 84856.164
$ws.Range("I103").Value = 333732.66
$ws.Range("J103").Value = 1897.3334
$ws.Range("K103").Value = 1001197.98
$ws.Range("L103").Value = 5692.0002
$ws.Range("M103").Value = -1000611.98
$ws.Range("N103").Value = -6864.0002
$ws.Range("H122").Value = 1365
$ws.Range("I122").Value = 1365
$ws.Range("K122").Value = 4095
$ws.Range("M122").Value = -1645
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2564.81
$ws.Range("I32").Value = 2564.81
$ws.Range("K32").Value = 2564.81
$ws.Range("M32").Value = -2277.81
$ws.Range("H61").Value = 20838418
$ws.Range("I61").Value = 83337000
$ws.Range("J61").Value = 5558.3335
$ws.Range("K61").Value = 83337000
$ws.Range("L61").Value = 5558.3335
$ws.Range("M61").Value = -83336788
$ws.Range("N61").Value = -5982.3335
$ws.Range("H136").Value = 20838418
$ws.Range("I136").Value = 83337000
$ws.Range("J136").Value = 5558.3335
$ws.Range("K136").Value = 250011000
$ws.Range("L136").Value = 16675.0005
$ws.Range("M136").Value = -250008450
$ws.Range("N136").Value = -21775.0005
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 463.1
$ws.Range("I94").Value = 463.1
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 463.1
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -12.10000000000002
$ws.Range("N94").ClearContents()
$ws.Range("H99").Value = 1364.5454
$ws.Range("I99").Value = 1402
$ws.Range("J99").Value = 1333.3334
$ws.Range("K99").Value = 1402
$ws.Range("L99").Value = 1333.3334
$ws.Range("M99").Value = 96
$ws.Range("N99").Value = -4329.3334
$ws.Range("H107").Value = 68626.92999999999
$ws.Range("I107").Value = 101931.4
$ws.Range("J107").Value = 2018
$ws.Range("K107").Value = 101931.4
$ws.Range("L107").Value = 2018
$ws.Range("M107").Value = -100011.4
$ws.Range("N107").Value = -5858
$ws.Range("H134").Value = 3008.9697
$ws.Range("I134").Value = 3063.8262
$ws.Range("J134").Value = 2882.8
$ws.Range("K134").Value = 9191.4786
$ws.Range("L134").Value = 8648.400000000001
$ws.Range("M134").Value = -6656.4786
$ws.Range("N134").Value = -13718.4
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6024.5244
$ws.Range("I31").Value = 2529.4138
$ws.Range("J31").Value = 7936.9434
$ws.Range("K31").Value = 2529.4138
$ws.Range("L31").Value = 7936.9434
$ws.Range("M31").Value = -2234.4138
$ws.Range("N31").Value = -8526.9434
$ws.Range("H34").Value = 6024.5244
$ws.Range("I34").Value = 2529.4138
$ws.Range("J34").Value = 7936.9434
$ws.Range("K34").Value = 2529.4138
$ws.Range("L34").Value = 7936.9434
$ws.Range("M34").Value = -2327.4138
$ws.Range("N34").Value = -8340.9434
$ws.Range("H58").Value = 1495.826
$ws.Range("I58").Value = 1103.6
$ws.Range("J58").Value = 2231.25
$ws.Range("K58").Value = 1103.6
$ws.Range("L58").Value = 2231.25
$ws.Range("M58").Value = -900.5999999999999
$ws.Range("N58").Value = -2637.25
$ws.Range("H107").Value = 2718080.8
$ws.Range("I107").Value = 3907026
$ws.Range("K107").Value = 3907026
$ws.Range("M107").Value = -3905106
$ws.Range("H124").Value = 30000
$ws.Range("I124").Value = 30000
$ws.Range("K124").Value = 30000
$ws.Range("M124").Value = -27545
$ws.Range("H127").Value = 88780
$ws.Range("J127").Value = 88780
$ws.Range("L127").Value = 88780
$ws.Range("N127").Value = -98700
$ws.Range("H130").Value = 88780
$ws.Range("J130").Value = 88780
$ws.Range("L130").Value = 88780
$ws.Range("N130").Value = -98820
$ws.Range("H136").Value = 1495.826
$ws.Range("I136").Value = 1103.6
$ws.Range("J136").Value = 2231.25
$ws.Range("K136").Value = 3310.8
$ws.Range("L136").Value = 6693.75
$ws.Range("M136").Value = -760.7999999999997
$ws.Range("N136").Value = -11793.75
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 59.652172
$ws.Range("I12").Value = 28.647058
$ws.Range("J12").Value = 147.5
$ws.Range("K12").Value = 85.941174
$ws.Range("L12").Value = 442.5
$ws.Range("M12").Value = 87.058826
$ws.Range("N12").Value = -788.5
$ws.Range("H34").Value = 7143403
$ws.Range("J34").Value = 8065104.5
$ws.Range("L34").Value = 24195313.5
$ws.Range("N34").Value = -24195481.5
$ws.Range("H39").Value = 2012.6
$ws.Range("J39").Value = 2012.6
$ws.Range("L39").Value = 6037.799999999999
$ws.Range("N39").Value = -6625.799999999999
$ws.Range("H55").Value = 1723.7059
$ws.Range("J55").Value = 1743.9375
$ws.Range("L55").Value = 5231.8125
$ws.Range("N55").Value = -5585.8125
$ws.Range("H113").Value = 512
$ws.Range("I113").Value = 484.97675
$ws.Range("J113").Value = 570.1
$ws.Range("K113").Value = 1454.93025
$ws.Range("L113").Value = 1710.3
$ws.Range("M113").Value = 715.0697500000001
$ws.Range("N113").Value = -6050.3
$ws.Range("H131").Value = 3809.8333
$ws.Range("I131").Value = 643.3333
$ws.Range("J131").Value = 4053.4102
$ws.Range("K131").Value = 1929.9999
$ws.Range("L131").Value = 12160.2306
$ws.Range("M131").Value = 3110.0001
$ws.Range("N131").Value = -22240.2306
$ws.Range("H132").Value = 2738.9473
$ws.Range("I132").Value = 2646.818
$ws.Range("J132").Value = 2865.625
$ws.Range("K132").Value = 23821.362
$ws.Range("L132").Value = 25790.625
$ws.Range("M132").Value = -21291.362
$ws.Range("N132").Value = -30850.625
$ws.Range("H134").Value = 5667.712
$ws.Range("I134").Value = 1933.6666
$ws.Range("J134").Value = 8818.3125
$ws.Range("K134").Value = 5800.9998
$ws.Range("L134").Value = 26454.9375
$ws.Range("M134").Value = -730.9997999999996
$ws.Range("N134").Value = -36594.9375
$ws.Range("H139").Value = 315990.34
$ws.Range("I139").Value = 436541.34
$ws.Range("J139").Value = 7915.5557
$ws.Range("K139").Value = 1309624.02
$ws.Range("L139").Value = 23746.6671
$ws.Range("M139").Value = -1304484.02
$ws.Range("N139").Value = -34026.6671
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2219.6
$ws.Range("I102").Value = 2184.1538
$ws.Range("J102").Value = 2450
$ws.Range("K102").Value = 2184.1538
$ws.Range("L102").Value = 2450
$ws.Range("M102").Value = -562.1538
$ws.Range("N102").Value = -5694
$ws.Range("H107").Value = 365.05554
$ws.Range("I107").Value = 267.07693
$ws.Range("J107").Value = 619.8
$ws.Range("K107").Value = 267.07693
$ws.Range("L107").Value = 619.8
$ws.Range("M107").Value = 1652.92307
$ws.Range("N107").Value = -4459.8
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8367.5
$ws.Range("I7").Value = 6333.3335
$ws.Range("J7").Value = 10401.667
$ws.Range("K7").Value = 6333.3335
$ws.Range("L7").Value = 10401.667
$ws.Range("M7").Value = -6221.3335
$ws.Range("N7").Value = -10625.667
$ws.Range("H109").Value = 29250
$ws.Range("J109").Value = 29250
$ws.Range("L109").Value = 29250
$ws.Range("N109").Value = -32024
$ws.Range("H122").Value = 5086.119
$ws.Range("I122").Value = 3739.8667
$ws.Range("J122").Value = 5834.037
$ws.Range("K122").Value = 11219.6001
$ws.Range("L122").Value = 17502.111
$ws.Range("M122").Value = -8769.6001
$ws.Range("N122").Value = -22402.111
$ws.Range("H126").Value = 8367.5
$ws.Range("I126").Value = 6333.3335
$ws.Range("J126").Value = 10401.667
$ws.Range("K126").Value = 19000.0005
$ws.Range("L126").Value = 31205.001
$ws.Range("M126").Value = -16530.0005
$ws.Range("N126").Value = -36145.001
$ws.Range("H132").Value = 2621.55
$ws.Range("I132").Value = 1764.0769
$ws.Range("J132").Value = 4214
$ws.Range("K132").Value = 5292.2307
$ws.Range("L132").Value = 12642
$ws.Range("M132").Value = -2762.2307
$ws.Range("N132").Value = -17702
$ws.Range("H136").Value = 6946957.5
$ws.Range("I136").Value = 2127.1333
$ws.Range("J136").Value = 18521674
$ws.Range("K136").Value = 6381.3999
$ws.Range("L136").Value = 55565022
$ws.Range("M136").Value = -3831.3999
$ws.Range("N136").Value = -55570122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 67204.2
$ws.Range("J15").Value = 67204.2
$ws.Range("L15").Value = 67204.2
$ws.Range("N15").Value = -67780.2
$ws.Range("H122").Value = 4300.4443
$ws.Range("I122").Value = 3751
$ws.Range("J122").Value = 4740
$ws.Range("K122").Value = 11253
$ws.Range("L122").Value = 14220
$ws.Range("M122").Value = -8803
$ws.Range("N122").Value = -19120
$ws.Range("H132").Value = 3943469.5
$ws.Range("I132").Value = 1898.6296
$ws.Range("K132").Value = 5695.8888
$ws.Range("M132").Value = -3165.8888
